$wb = $excel.ActiveWorkbook

# Helper: write a literal value into a cell while keeping it stored as TEXT
# (matches the workbook's existing convention of storing every value -
# numeric-looking or not - as a shared string). Temporarily switching the
# cell to a Text number format stops Excel's "smart" numeric auto-detection
# from converting e.g. "-8.8" into a real number; resetting the style back
# to "Normal" afterwards drops the temporary format so the cell ends up with
# its original (default) style, same as before the edit.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Sheet: Restricciones_del_follower (index 3)
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("A2") "11.3 - 2x_1 + y_1 - y_2"
Set-TextValue $ws3.Range("B2") "-8.8"
Set-TextValue $ws3.Range("D2") "0.79"
Set-TextValue $ws3.Range("F2") "3.5"
Set-TextValue $ws3.Range("A3") "-3.55 + x_1 - 3x_2 + y_2"
Set-TextValue $ws3.Range("B3") "1.5499999999999998"
Set-TextValue $ws3.Range("D3") "0.09"
Set-TextValue $ws3.Range("E3") "8.100000000000001"
Set-TextValue $ws3.Range("A4") "-9.129999999999999 + x_1 + x_2"
Set-TextValue $ws3.Range("B4") "6.75"
Set-TextValue $ws3.Range("D4") "0.54"
Set-TextValue $ws3.Range("E4") "0.4"
Set-TextValue $ws3.Range("F4") "3.8"

# Sheet: Punto_modificado (index 4)
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "6.65"
Set-TextValue $ws4.Range("B2") "2.1"
Set-TextValue $ws4.Range("C2") "5.2"
Set-TextValue $ws4.Range("D2") "3.2"

# Sheet: Vector_bf (index 5; name collides case-insensitively with Vector_BF)
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "3.21"
Set-TextValue $ws5.Range("A3") "-0.29999999999999993"

# Sheet: Vector_BF (index 6)
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-6.500000000000002"
Set-TextValue $ws6.Range("A3") "22.900000000000006"
Set-TextValue $ws6.Range("A5") "-8.100000000000001"
